$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the two new hidden helper cells (row 4) with the additional
# data-binding placeholders for course code and class code.
$ws.Range("F4").Value = "&=[DATA].MAKHOAHOC"
$ws.Range("I4").Value = "&=[DATA].MALOPHOC"

# Reflect the cursor/selection position that was active when the sheet
# was last saved (also clears the previous "scrolled to column H" view).
$ws.Range("H15").Select()
